# Add the new "High School Units" sheet (sheetId=3) at the end of the workbook,
# populate its headers/data, size a couple of columns, and select it as active -
# mirroring the commit "finished high school units and test".

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$ws.Name = "High School Units"

# Header row
$headers = @(
    "Value",
    "english",
    "mathematics",
    "science",
    "lab",
    "foreign-language",
    "social-studies",
    "history",
    "academic-electives",
    "computer-science",
    "visual/performing-arts",
    "others",
    "units-required",
    "units-recommended"
)
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows (rows 2-24, columns A-N)
$data = @(
    @(4,1,0,0,0,0,0,0,0,0,0,0,1,0),
    @(0,1,0,0,0,0,0,0,0,0,0,0,0,1),
    @(4,0,1,0,0,0,0,0,0,0,0,0,1,0),
    @(5,0,1,0,0,0,0,0,0,0,0,0,0,1),
    @(3,0,0,1,0,0,0,0,0,0,0,0,1,0),
    @(4,0,1,0,0,0,0,0,0,0,0,0,0,1),
    @(3,0,0,1,1,0,0,0,0,0,0,0,1,0),
    @(0,0,0,0,0,1,0,0,0,0,0,0,1,0),
    @(0,0,0,0,0,1,0,0,0,0,0,0,0,1),
    @(2,0,0,0,0,0,1,0,0,0,0,0,1,0),
    @(0,0,0,0,0,0,1,0,0,0,0,0,0,1),
    @(0,0,0,0,0,0,0,1,0,0,0,0,1,0),
    @(0,0,0,0,0,0,0,1,0,0,0,0,0,1),
    @(0,0,0,0,0,0,0,1,0,0,0,0,1,0),
    @(0,0,0,0,0,0,0,1,0,0,0,0,0,1),
    @(0,0,0,0,0,0,0,0,1,0,0,0,1,0),
    @(4,0,0,0,0,0,0,0,1,0,0,0,0,1),
    @(0,0,0,0,0,0,0,0,0,1,0,0,1,0),
    @(0,0,0,0,0,0,0,0,0,1,0,0,0,1),
    @(0,0,0,0,0,0,0,0,0,0,1,0,1,0),
    @(0,0,0,0,0,0,0,0,0,0,1,0,0,1),
    @(0,0,0,0,0,0,0,0,0,0,0,1,1,0),
    @(0,0,0,0,0,0,0,0,0,0,0,1,0,1)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}

# Column sizing to match the authored widths (M ~17 chars, N ~21.33 chars)
$ws.Columns.Item(13).ColumnWidth = 16.166666666666668
$ws.Columns.Item(14).ColumnWidth = 20.498697916666668

# Make the new sheet the active / selected tab with A19 selected
[void]$ws.Range("A19").Select()
